$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $scoped = $d.Paragraphs.Item($paraIndex).Range.Duplicate()
    $ok = $scoped.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for paragraph $paraIndex"
    }
}

# --- Single-run whole-paragraph replacements ---
Replace-InParagraph 6 "Introduzir o aluno ao conhecimento e entendimento das funções da química orgânica. Compreender as condições para que as reações orgânicas ocorram. Saber analisar os produtos das reações empregando espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono. Levá-los a compreender os procedimentos e problemas industriais, sociais e ambientais, com os quais a Engenharia Química está estritamente relacionada, tornando-os, dessa forma, capazes a exercerem a função do Engenheiro Químico e realizar as mudanças que se fizerem necessárias." "Ácidos carboxílicos e seus derivados (haletos de acila, anidridos, ésteres e amidas – 16 horas/aulas). Aldeídos e cetonas (4 horas/aulas). Aminas (4 horas/aulas). Fenóis (3 horas/aulas). Sais de diazônio e azo compostos (3 horas/aulas). Ácidos sulfônicos e derivados (4 horas/aulas). Cetoácidos (3 horas/aulas). Hidroxi ácidos (3 aulas). Heterocíclicos (4 horas/aulas). Noções de síntese orgânica (8 horas/aulas). Análises por espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono (8 horas/aulas)"
Replace-InParagraph 7 "Introduce the student to the knowledge and understanding of the functions of organic chemistry. Understand the conditions for organic reactions to occur. Know how to analyze the reaction products using spectroscopy in the infrared region and proton and carbon nuclear magnetic resonance. Get them to understand the industrial, social and environmental procedures and problems with which Chemical Engineering is closely related, thus making them able to perform the function of the Chemical Engineer and make the necessary changes." "Carboxylic acids and their derivatives (acyl halides, anhydrides, esters and amides) - 16 hours. Aldehydes and Ketones (4 hours). Amines (4 hours). Phenols (3 hours). Diazonium salts and azo compounds (3 hours). Sulfonic acids and derivatives (4 hours). Keto acids (3 hours). Hydroxy acids (3 hours). Heterocyclic (4 hours). Notion for organic synthesis (8 hours). Analyzes by infrared spectroscopy and proton and carbon nuclear magnetic resonance (8 hours)."
Replace-InParagraph 9 "5840751 - Jayne Carlos de Souza Barboza" "Introduzir o aluno ao conhecimento e entendimento das funções da química orgânica. Compreender as condições para que as reações orgânicas ocorram. Saber analisar os produtos das reações empregando espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono. Levá-los a compreender os procedimentos e problemas industriais, sociais e ambientais, com os quais a Engenharia Química está estritamente relacionada, tornando-os, dessa forma, capazes a exercerem a função do Engenheiro Químico e realizar as mudanças que se fizerem necessárias."
Replace-InParagraph 11 "Ácidos carboxílicos e seus derivados (haletos de acila, anidridos, ésteres e amidas – 16 horas/aulas). Aldeídos e cetonas (4 horas/aulas). Aminas (4 horas/aulas). Fenóis (3 horas/aulas). Sais de diazônio e azo compostos (3 horas/aulas). Ácidos sulfônicos e derivados (4 horas/aulas). Cetoácidos (3 horas/aulas). Hidroxi ácidos (3 aulas). Heterocíclicos (4 horas/aulas). Noções de síntese orgânica (8 horas/aulas). Análises por espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono (8 horas/aulas)" "Ácidos carboxílicos e derivados: Nomenclatura, propriedades físicas, processos de obtenção e propriedades químicas. Reações de substituição nucleofilícas em carbono acíclico. Síntese dos haletos de acila, anidridos, ésteres, amidas, aldeídos, cetonas, aminas, fenóis e suas propriedades químicas. Alfa halo-ácidos. Cetoácidos: processos de obtenção e propriedades químicas. Preparação de alfa cetoésteres. Hidroxi ácidos: processos de obtenção e propriedades químicas. Ácidos sulfônicos e seus derivados: processos de obtenção e propriedades químicas. Cloreto de sulfonila, sulfonamidas. Ésteres dos ácidos sulfônicos. Aminas: propriedades físicas, processos de obtenção e propriedades químicas. Reações de acoplamento dos sais de diazônio e azo compostos: Fenóis. Noções de síntese orgânica. Preparação de amostras, análises nos equipamentos e análises de espectros de compostos orgânicos na região do infravermelho e por ressonância magnética nuclear de próton e carbono"
Replace-InParagraph 12 "Carboxylic acids and their derivatives (acyl halides, anhydrides, esters and amides) - 16 hours. Aldehydes and Ketones (4 hours). Amines (4 hours). Phenols (3 hours). Diazonium salts and azo compounds (3 hours). Sulfonic acids and derivatives (4 hours). Keto acids (3 hours). Hydroxy acids (3 hours). Heterocyclic (4 hours). Notion for organic synthesis (8 hours). Analyzes by infrared spectroscopy and proton and carbon nuclear magnetic resonance (8 hours)." "Introduce the student to the knowledge and understanding of the functions of organic chemistry. Understand the conditions for organic reactions to occur. Know how to analyze the reaction products using spectroscopy in the infrared region and proton and carbon nuclear magnetic resonance. Get them to understand the industrial, social and environmental procedures and problems with which Chemical Engineering is closely related, thus making them able to perform the function of the Chemical Engineer and make the necessary changes."
Replace-InParagraph 19 "1) DURST, H. D. Fundamentos de Química Orgânica. São Paulo: Editora Edgard Blucher, 1997. 2) HENDRIKSON, J. B.; CRAM, D. J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966. 3) SOLOMONS, T.W.G; FRYHLE, G.. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001. 4)VOLLHARDT, K. PETER C. Química Orgânica: estrutura e função.Bookman, Porto Alegre 4 edição, 2004. 5) Silverstein, R. M.; Bassler, G. C.; Morrill, T. C. Identificação Espectrométrica de Compostos Orgânicos. Guanabara Koogan, 5 edição, 1991.6) PAVIA, D,L,; LAMPMAN, G.M.; KRIZ, G.S.  Introduction to Spectroscopy- A Guide for Students of Organic Chemistry. Harcourt Brace College Publisheres, 1996" "5840751 - Jayne Carlos de Souza Barboza"

# --- Word paragraph 14 (Programa resumido body text) ---
Replace-InParagraph 14 "Ácidos carboxílicos e derivados: Nomenclatura, propriedades físicas, processos de obtenção e propriedades químicas. Reações de substituição nucleofilícas em carbono acíclico. Síntese dos haletos de acila, anidridos, ésteres, amidas, aldeídos, cetonas, aminas, fenóis e suas propriedades químicas. Alfa halo-ácidos. Cetoácidos: processos de obtenção e propriedades químicas. Preparação de alfa cetoésteres. Hidroxi ácidos: processos de obtenção e propriedades químicas. Ácidos sulfônicos e seus derivados: processos de obtenção e propriedades químicas. Cloreto de sulfonila, sulfonamidas. Ésteres dos ácidos sulfônicos. Aminas: propriedades físicas, processos de obtenção e propriedades químicas. Reações de acoplamento dos sais de diazônio e azo compostos: Fenóis. Noções de síntese orgânica. Preparação de amostras, análises nos equipamentos e análises de espectros de compostos orgânicos na região do infravermelho e por ressonância magnética nuclear de próton e carbono" "Duas provas semestrais teóricas (P1 e P2)."

# --- Word paragraph 17 (Avaliacao bullet: 3 value-runs) ---
Replace-InParagraph 17 "Duas provas semestrais teóricas (P1 e P2)." "A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3"
Replace-InParagraph 17 "A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3" "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
Replace-InParagraph 17 "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno." "1) DURST, H. D. Fundamentos de Química Orgânica. São Paulo: Editora Edgard Blucher, 1997. 2) HENDRIKSON, J. B.; CRAM, D. J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966. 3) SOLOMONS, T.W.G; FRYHLE, G.. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001. 4)VOLLHARDT, K. PETER C. Química Orgânica: estrutura e função.Bookman, Porto Alegre 4 edição, 2004. 5) Silverstein, R. M.; Bassler, G. C.; Morrill, T. C. Identificação Espectrométrica de Compostos Orgânicos. Guanabara Koogan, 5 edição, 1991.6) PAVIA, D,L,; LAMPMAN, G.M.; KRIZ, G.S.  Introduction to Spectroscopy- A Guide for Students of Organic Chemistry. Harcourt Brace College Publisheres, 1996"

Write-Output "Done"
